$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.149.43"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.637.75"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "216.67"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "0.519"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.865.20"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.645.87"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "0.541"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "66.53"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "27.136.52"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").Value = "217.01"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").Value = "146.69"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "7.40"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "1.300.76"
$ws.Range("E34").Value = "  +2.82%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "2.47"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "0.854"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").Value = "0.544"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  +5.80%  "
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "1.776.58"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "61.65"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").Value = "91.39"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "7.66"
$ws.Range("E50").Value = "  -0.24%  "
